# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2-8, columns B (TB), C (d2S), D (K), E (IP); F (Win) unchanged; G (sum) recalculated as B+C+D+E
$data = @(
    @{Row=2; B=3.182878228561681;  C=1.65323645889881;   D=0.1529057820181812; E=0.4998867070740569; G=5.488907176552729}
    @{Row=3; B=3.182878228561681;  C=1.65323645889881;   D=0.7127328510149897; E=0.4998867070740569; G=6.048734245549538}
    @{Row=4; B=0.02258322285507441;C=0.3375848360084654; D=0.1529057820181812; E=0.4998867070740569; G=1.012960547955778}
    @{Row=5; B=1.505614041169197;  C=1.65323645889881;   D=3.082599426703578;  E=0.4998867070740569; G=6.741336633845642}
    @{Row=6; B=3.182878228561681;  C=1.65323645889881;   D=0.7127328510149897; E=0.4998867070740569; G=6.048734245549538}
    @{Row=7; B=3.182878228561681;  C=0.05231270169004087;D=0.1529057820181812; E=0.4998867070740569; G=3.887983419343961}
    @{Row=8; B=3.182878228561681;  C=1.65323645889881;   D=0.7127328510149897; E=0.4998867070740569; G=6.048734245549538}
)

foreach ($rowData in $data) {
    $r = $rowData.Row
    $ws.Range("B$r").Value = $rowData.B
    $ws.Range("C$r").Value = $rowData.C
    $ws.Range("D$r").Value = $rowData.D
    $ws.Range("E$r").Value = $rowData.E
    $ws.Range("G$r").Value = $rowData.G
}
